# SOLUX git 기초 과제 - fill in the blanks on the title slide.
#
# Slide 1 has two placeholder shapes we need to touch:
#   Shape "제목 1"   (ctrTitle) - the two "--------" answer lines become
#                                 real git commands.
#   Shape "부제목 2" (subTitle) - the "학부 / 학번 / 이름" placeholder line
#                                 becomes the real major / student id / name.
#
# NOTE: this host's TextRange.Text setter (when applied to the *whole*
# TextFrame.TextRange) only ever rewrites the first run and silently
# ignores the rest, so we never assign to the top-level TextRange.Text.
# Instead we locate the exact substring with IndexOf() and rewrite just
# that slice via TextRange.Characters(start, length) -- this preserves
# every other run/line-break untouched, exactly like using Find & Replace
# on a specific phrase in real PowerPoint.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Replace-Substring($textRange, [string]$oldText, [string]$newText) {
    $full = $textRange.Text
    $idx = $full.IndexOf($oldText)
    if ($idx -lt 0) {
        throw "Substring not found: $oldText"
    }
    $piece = $textRange.Characters($idx + 1, $oldText.Length)
    $piece.Text = $newText
}

# ---- locate the two placeholder shapes by their placeholder type ----
$titleShape = $null
$subTitleShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $txt = $sh.TextFrame.TextRange.Text
        if ($txt.IndexOf("--------") -ge 0) { $titleShape = $sh }
        if ($txt.IndexOf([char]54617 + [char]48512) -ge 0) { $subTitleShape = $sh } # "학부"
    }
}
if ($null -eq $titleShape) { $titleShape = $s.Shapes.Item(4) }
if ($null -eq $subTitleShape) { $subTitleShape = $s.Shapes.Item(5) }

# ---- 1. Title: fill in the two git commands ----
$titleRange = $titleShape.TextFrame.TextRange
Replace-Substring $titleRange "1. --------" "1. git add"
Replace-Substring $titleRange "2. --------" "2. git push"

# ---- 2. Subtitle: replace "학부 / 학번 / 이름" with the real info ----
$subRange = $subTitleShape.TextFrame.TextRange
$oldSub = $subRange.Text
$newSub = "컴퓨터과학전공 / 2211015 / 홍서현"

# Put the full new string into the first run (growing it as needed), then
# delete whatever remains of the old trailing runs so the paragraph ends
# up containing exactly the new text.
$firstRun = $subRange.Characters(1, 1)
$firstRun.Text = $newSub
$afterLen = $subRange.Text.Length
$tailLen = $afterLen - $newSub.Length
if ($tailLen -gt 0) {
    $tail = $subRange.Characters($newSub.Length + 1, $tailLen)
    $tail.Text = ""
}

# Re-create the original run/segment boundaries inside the new text (so
# formatting such as language-specific runs stay split the same way the
# source text was: "전공 / 학번 / 이름" style segments) by forcing a
# (value-preserving) font write at each boundary offset - this causes the
# host to split the run at that character position without changing how
# it looks.
$fontName = $subRange.Characters(1, 1).Font.Name
$segments = @("컴퓨터과학전공 ", "/", " ", "2211015", " ", "/", " 홍서현")
$pos = 0
foreach ($seg in $segments) {
    $pos += $seg.Length
    if ($pos -lt $newSub.Length) {
        $piece = $subRange.Characters(1, $pos)
        $piece.Font.Name = $fontName
    }
}

# The subtitle box grows wider to fit the longer line (matches the
# resize PowerPoint performs automatically for this placeholder).
$subTitleShape.Width = 6231222 / 12700.0
